# Case_3_177: update pl_mw.xlsx result values for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Build a 24x11 array covering B2:L25 (columns B..L, rows 2..25)
$arr = New-Object 'object[,]' 24,11

# Row 2 (A2 index 0)
$arr[0,0] = 2.073190472657984
$arr[0,1] = 0.1611398126446772
$arr[0,2] = 0.09127684349529019
$arr[0,3] = 0
$arr[0,4] = 2.345165255646435
$arr[0,5] = 1.676800007928222
$arr[0,6] = 1.473155296510186
$arr[0,7] = 0
$arr[0,8] = 0.1913408601438089
$arr[0,9] = 0
$arr[0,10] = 0.3872427179576903
# Row 3 (A3 index 1)
$arr[1,0] = 1.95656535736407
$arr[1,1] = 0.1441581447731153
$arr[1,2] = 0.09090038879359597
$arr[1,3] = 0
$arr[1,4] = 2.347493863243159
$arr[1,5] = 1.674475025908563
$arr[1,6] = 1.479225310039354
$arr[1,7] = 0
$arr[1,8] = 0.1930563297397931
$arr[1,9] = 0
$arr[1,10] = 0.3798261718278297
# Row 4 (A4 index 2)
$arr[2,0] = 1.885725325509839
$arr[2,1] = 0.133672008254365
$arr[2,2] = 0.09067935948151984
$arr[2,3] = 0
$arr[2,4] = 2.350371484442704
$arr[2,5] = 1.674265717346486
$arr[2,6] = 1.483775618622062
$arr[2,7] = 0
$arr[2,8] = 0.1941777628944834
$arr[2,9] = 0
$arr[2,10] = 0.3754265323428854
# Row 5 (A5 index 3)
$arr[3,0] = 1.857051517381478
$arr[3,1] = 0.1293839957401701
$arr[3,2] = 0.09059184975274093
$arr[3,3] = 0
$arr[3,4] = 2.351907474466913
$arr[3,5] = 1.674485754647648
$arr[3,6] = 1.485836617885639
$arr[3,7] = 0
$arr[3,8] = 0.1946518946041724
$arr[3,9] = 0
$arr[3,10] = 0.37367250810928
# Row 6 (A6 index 4)
$arr[4,0] = 1.852302002411648
$arr[4,1] = 0.1286710817993253
$arr[4,2] = 0.09057747407240413
$arr[4,3] = 0
$arr[4,4] = 2.352184447540779
$arr[4,5] = 1.674540703639309
$arr[4,6] = 1.486191320991537
$arr[4,7] = 0
$arr[4,8] = 0.1947316591574335
$arr[4,9] = 0
$arr[4,10] = 0.3733836041171514
# Row 7 (A7 index 5)
$arr[5,0] = 1.885337833478559
$arr[5,1] = 0.1336142385247285
$arr[5,2] = 0.09067816889870883
$arr[5,3] = 0
$arr[5,4] = 2.350390729196576
$arr[5,5] = 1.674267449907205
$arr[5,6] = 1.483802577436961
$arr[5,7] = 0
$arr[5,8] = 0.1941840877974688
$arr[5,9] = 0
$arr[5,10] = 0.3754027194589753
# Row 8 (A8 index 6)
$arr[6,0] = 2.032819167844821
$arr[6,1] = 0.1552968944695863
$arr[6,2] = 0.09114495291554903
$arr[6,3] = 0
$arr[6,4] = 2.34566720914745
$arr[6,5] = 1.675744898377218
$arr[6,6] = 1.475077171409893
$arr[6,7] = 0
$arr[6,8] = 0.1919182186128623
$arr[6,9] = 0
$arr[6,10] = 0.3846535565816822
# Row 9 (A9 index 7)
$arr[7,0] = 2.32810404988254
$arr[7,1] = 0.1973440037238561
$arr[7,2] = 0.09213987700297466
$arr[7,3] = 0
$arr[7,4] = 2.347928519379877
$arr[7,5] = 1.68835626697981
$arr[7,6] = 1.464514404433629
$arr[7,7] = 0
$arr[7,8] = 0.188015058170647
$arr[7,9] = 0
$arr[7,10] = 0.4040144515012969
# Row 10 (A10 index 8)
$arr[8,0] = 2.548750816030122
$arr[8,1] = 0.2279487546029202
$arr[8,2] = 0.09291854064959182
$arr[8,3] = 0
$arr[8,4] = 2.356669461281967
$arr[8,5] = 1.70361393355725
$arr[8,6] = 1.460768367579419
$arr[8,7] = 0
$arr[8,8] = 0.1854761800365274
$arr[8,9] = 0
$arr[8,10] = 0.4189806271361647
# Row 11 (A11 index 9)
$arr[9,0] = 2.649934126848962
$arr[9,1] = 0.241809711273504
$arr[9,2] = 0.09328297395022744
$arr[9,3] = 0
$arr[9,4] = 2.362195195546477
$arr[9,5] = 1.711871546500248
$arr[9,6] = 1.459940950688036
$arr[9,7] = 0
$arr[9,8] = 0.1843924534483552
$arr[9,9] = 0
$arr[9,10] = 0.4259499342398243
# Row 12 (A12 index 10)
$arr[10,0] = 2.68836572958952
$arr[10,1] = 0.2470496538507234
$arr[10,2] = 0.09342242966280523
$arr[10,3] = 0
$arr[10,4] = 2.364511362379801
$arr[10,5] = 1.715189028635706
$arr[10,6] = 1.459754085665736
$arr[10,7] = 0
$arr[10,8] = 0.1839923094945881
$arr[10,9] = 0
$arr[10,10] = 0.4286121391527615
# Row 13 (A13 index 11)
$arr[11,0] = 2.680083663777793
$arr[11,1] = 0.2459215356579989
$arr[11,2] = 0.09339233102309663
$arr[11,3] = 0
$arr[11,4] = 2.364002571549733
$arr[11,5] = 1.714466058524806
$arr[11,6] = 1.459788699732115
$arr[11,7] = 0
$arr[11,8] = 0.1840780322643063
$arr[11,9] = 0
$arr[11,10] = 0.4280377606104366
# Row 14 (A14 index 12)
$arr[12,0] = 2.653093604000333
$arr[12,1] = 0.2422409833686459
$arr[12,2] = 0.0932944180313342
$arr[12,3] = 0
$arr[12,4] = 2.362381259800003
$arr[12,5] = 1.712140653205495
$arr[12,6] = 1.459923040474564
$arr[12,7] = 0
$arr[12,8] = 0.1843593281533877
$arr[12,9] = 0
$arr[12,10] = 0.4261684933916854
# Row 15 (A15 index 13)
$arr[13,0] = 2.63657645481203
$arr[13,1] = 0.2399853771734399
$arr[13,2] = 0.09323463221182493
$arr[13,3] = 0
$arr[13,4] = 2.361417317057814
$arr[13,5] = 1.710741118528915
$arr[13,6] = 1.460021808632405
$arr[13,7] = 0
$arr[13,8] = 0.1845329634787198
$arr[13,9] = 0
$arr[13,10] = 0.4250265162912683
# Row 16 (A16 index 14)
$arr[14,0] = 2.542154478893167
$arr[14,1] = 0.2270416675325748
$arr[14,2] = 0.09289492823839751
$arr[14,3] = 0
$arr[14,4] = 2.35633959563782
$arr[14,5] = 1.703100874837361
$arr[14,6] = 1.460840115853443
$arr[14,7] = 0
$arr[14,8] = 0.1855484371227707
$arr[14,9] = 0
$arr[14,10] = 0.4185284009649308
# Row 17 (A17 index 15)
$arr[15,0] = 2.484436535182056
$arr[15,1] = 0.2190853581870158
$arr[15,2] = 0.09268913556207536
$arr[15,3] = 0
$arr[15,4] = 2.353622053312264
$arr[15,5] = 1.69875191687666
$arr[15,6] = 1.461566938305083
$arr[15,7] = 0
$arr[15,8] = 0.1861896375767795
$arr[15,9] = 0
$arr[15,10] = 0.4145832197470867
# Row 18 (A18 index 16)
$arr[16,0] = 2.451314986125396
$arr[16,1] = 0.2145033346999696
$arr[16,2] = 0.09257173151198828
$arr[16,3] = 0
$arr[16,4] = 2.352204781986444
$arr[16,5] = 1.696374388294061
$arr[16,6] = 1.46206748029843
$arr[16,7] = 0
$arr[16,8] = 0.1865651429536115
$arr[16,9] = 0
$arr[16,10] = 0.4123292262643758
# Row 19 (A19 index 17)
$arr[17,0] = 2.44011372567843
$arr[17,1] = 0.2129509528753601
$arr[17,2] = 0.09253214632217066
$arr[17,3] = 0
$arr[17,4] = 2.35174993046931
$arr[17,5] = 1.695590639106001
$arr[17,6] = 1.462251111454293
$arr[17,7] = 0
$arr[17,8] = 0.1866934340695465
$arr[17,9] = 0
$arr[17,10] = 0.411568671094301
# Row 20 (A20 index 18)
$arr[18,0] = 2.49057281944016
$arr[18,1] = 0.2199329181669896
$arr[18,2] = 0.09271094305097449
$arr[18,3] = 0
$arr[18,4] = 2.353896245566759
$arr[18,5] = 1.699202041819746
$arr[18,6] = 1.461481026284446
$arr[18,7] = 0
$arr[18,8] = 0.1861206868875449
$arr[18,9] = 0
$arr[18,10] = 0.4150016215965024
# Row 21 (A21 index 19)
$arr[19,0] = 2.661018105492587
$arr[19,1] = 0.2433222930271768
$arr[19,2] = 0.09332313815150428
$arr[19,3] = 0
$arr[19,4] = 2.362851400375817
$arr[19,5] = 1.712818501759955
$arr[19,6] = 1.459880146095344
$arr[19,7] = 0
$arr[19,8] = 0.1842764268409098
$arr[19,9] = 0
$arr[19,10] = 0.4267169168088571
# Row 22 (A22 index 20)
$arr[20,0] = 2.77308787281811
$arr[20,1] = 0.2585566835211068
$arr[20,2] = 0.09373170345183723
$arr[20,3] = 0
$arr[20,4] = 2.370008369064308
$arr[20,5] = 1.722828579306679
$arr[20,6] = 1.459571168314682
$arr[20,7] = 0
$arr[20,8] = 0.1831307808626654
$arr[20,9] = 0
$arr[20,10] = 0.4345080255329208
# Row 23 (A23 index 21)
$arr[21,0] = 2.713212719473063
$arr[21,1] = 0.2504305804637568
$arr[21,2] = 0.09351287557823795
$arr[21,3] = 0
$arr[21,4] = 2.366068927090211
$arr[21,5] = 1.717383982376873
$arr[21,6] = 1.459668482150249
$arr[21,7] = 0
$arr[21,8] = 0.1837367730995005
$arr[21,9] = 0
$arr[21,10] = 0.430337487390787
# Row 24 (A24 index 22)
$arr[22,0] = 2.487798414433826
$arr[22,1] = 0.2195497607013124
$arr[22,2] = 0.09270108105359753
$arr[22,3] = 0
$arr[22,4] = 2.353771831404089
$arr[22,5] = 1.69899815808094
$arr[22,6] = 1.461519609587867
$arr[22,7] = 0
$arr[22,8] = 0.1861518380923339
$arr[22,9] = 0
$arr[22,10] = 0.4148124180596966
# Row 25 (A25 index 23)
$arr[23,0] = 2.247571709106694
$arr[23,1] = 0.1860196670606058
$arr[23,2] = 0.09186229139888269
$arr[23,3] = 0
$arr[23,4] = 2.346077457336193
$arr[23,5] = 1.683897360658904
$arr[23,6] = 1.466668590543236
$arr[23,7] = 0
$arr[23,8] = 0.1890131913223598
$arr[23,9] = 0
$arr[23,10] = 0.398646425029753

$ws.Range("B2:L25").Value = $arr
